# Added test cases in excel sheet:
#  - Projects sheet: point the ProjectUrl cell (C2) at the new SafeWay URL
#  - Widen column C on the Projects sheet so the longer URL fits
#  - Make "Projects" the active/selected sheet (was "SafeWay"), with E2 selected

$wb = $excel.ActiveWorkbook

$projects = $wb.Worksheets.Item("Projects")

# Update the project URL used for testing
$projects.Range("C2").Value = "http://safeway.com"

# Widen column C to fit the new (longer) URL text
$projects.Columns.Item(3).ColumnWidth = 46.3

# Make "Projects" the active sheet/tab and select cell E2 on it
$projects.Activate()
$projects.Range("E2").Select()
